$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6489.5713
$ws.Range("J17").Value = 6489.5713
$ws.Range("L17").Value = 19468.7139
$ws.Range("N17").Value = -19804.7139
$ws.Range("H18").Value = 805.1429000000001
$ws.Range("I18").Value = 805.1429000000001
$ws.Range("K18").Value = 805.1429000000001
$ws.Range("M18").Value = -521.1429000000001
$ws.Range("H98").Value = 828.10345
$ws.Range("I98").Value = 750.5714
$ws.Range("K98").Value = 750.5714
$ws.Range("M98").Value = 747.4286
$ws.Range("H111").Value = 256497.25
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 256497.25
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 769491.75
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -775625.75
$ws.Range("H122").Value = 828.10345
$ws.Range("I122").Value = 750.5714
$ws.Range("K122").Value = 2251.7142
$ws.Range("M122").Value = 198.2857999999997
$ws.Range("H127").Value = 2405.5715
$ws.Range("I127").Value = 1017.9
$ws.Range("K127").Value = 3053.7
$ws.Range("M127").Value = 1906.3
$ws.Range("H129").Value = 989
$ws.Range("I129").Value = 679.61536
$ws.Range("K129").Value = 2038.84608
$ws.Range("M129").Value = 2961.15392
$ws.Range("H131").Value = 2220.5833
$ws.Range("J131").Value = 3925
$ws.Range("L131").Value = 11775
$ws.Range("N131").Value = -21855
$ws.Range("H132").Value = 2106.889
$ws.Range("I132").Value = 2111.907
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 6335.721
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -3805.721
$ws.Range("N132").Value = -11057
$ws.Range("H133").Value = 124994.664
$ws.Range("J133").Value = 124994.664
$ws.Range("L133").Value = 124994.664
$ws.Range("N133").Value = -135114.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5458.2246
$ws.Range("J61").Value = 4571.4736
$ws.Range("L61").Value = 4571.4736
$ws.Range("N61").Value = -4995.4736
$ws.Range("H74").Value = 4776.0835
$ws.Range("I74").Value = 3092.85
$ws.Range("K74").Value = 3092.85
$ws.Range("M74").Value = -2218.85
$ws.Range("H77").Value = 4776.0835
$ws.Range("I77").Value = 3092.85
$ws.Range("K77").Value = 15464.25
$ws.Range("M77").Value = -11096.25
$ws.Range("H122").Value = 3805
$ws.Range("I122").Value = 2458.4285
$ws.Range("K122").Value = 7375.2855
$ws.Range("M122").Value = -4925.2855
$ws.Range("H133").Value = 91062.375
$ws.Range("J133").Value = 91062.375
$ws.Range("L133").Value = 91062.375
$ws.Range("N133").Value = -96122.375
$ws.Range("H136").Value = 5458.2246
$ws.Range("J136").Value = 4571.4736
$ws.Range("L136").Value = 13714.4208
$ws.Range("N136").Value = -18814.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1217.3158
$ws.Range("I99").Value = 1317.2858
$ws.Range("J99").Value = 937.4
$ws.Range("K99").Value = 1317.2858
$ws.Range("L99").Value = 937.4
$ws.Range("M99").Value = 180.7141999999999
$ws.Range("N99").Value = -3933.4
$ws.Range("H132").Value = 110397.336
$ws.Range("J132").Value = 110397.336
$ws.Range("L132").Value = 110397.336
$ws.Range("N132").Value = -120517.336
$ws.Range("H134").Value = 3157.5557
$ws.Range("I134").Value = 2412.7693
$ws.Range("J134").Value = 7998.6665
$ws.Range("K134").Value = 7238.3079
$ws.Range("L134").Value = 23995.9995
$ws.Range("M134").Value = -4703.3079
$ws.Range("N134").Value = -29065.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2227.889
$ws.Range("I22").Value = 410.6
$ws.Range("J22").Value = 4499.5
$ws.Range("K22").Value = 410.6
$ws.Range("L22").Value = 4499.5
$ws.Range("M22").Value = -60.60000000000002
$ws.Range("N22").Value = -5199.5
$ws.Range("H86").Value = 5355.231
$ws.Range("I86").Value = 4986.8335
$ws.Range("J86").Value = 5671
$ws.Range("K86").Value = 4986.8335
$ws.Range("L86").Value = 5671
$ws.Range("M86").Value = -3863.8335
$ws.Range("N86").Value = -7917
$ws.Range("H89").Value = 5355.231
$ws.Range("I89").Value = 4986.8335
$ws.Range("J89").Value = 5671
$ws.Range("K89").Value = 24934.1675
$ws.Range("L89").Value = 28355
$ws.Range("M89").Value = -19318.1675
$ws.Range("N89").Value = -39587
$ws.Range("H99").Value = 4308.686
$ws.Range("J99").Value = 2941.4348
$ws.Range("L99").Value = 2941.4348
$ws.Range("N99").Value = -5937.4348
$ws.Range("H122").Value = 2787.9285
$ws.Range("I122").Value = 2440.2
$ws.Range("J122").Value = 3189.1538
$ws.Range("K122").Value = 7320.599999999999
$ws.Range("L122").Value = 9567.4614
$ws.Range("M122").Value = -4870.599999999999
$ws.Range("N122").Value = -14467.4614
$ws.Range("H126").Value = 4308.686
$ws.Range("J126").Value = 2941.4348
$ws.Range("L126").Value = 8824.304400000001
$ws.Range("N126").Value = -13764.3044
$ws.Range("H132").Value = 4675.9
$ws.Range("I132").Value = 3169.025
$ws.Range("J132").Value = 10703.4
$ws.Range("K132").Value = 9507.075000000001
$ws.Range("L132").Value = 32110.2
$ws.Range("M132").Value = -6977.075000000001
$ws.Range("N132").Value = -37170.2
$ws.Range("H140").Value = 52146.6
$ws.Range("J140").Value = 52146.6
$ws.Range("L140").Value = 52146.6
$ws.Range("N140").Value = -62506.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2068.5715
$ws.Range("I138").Value = 2068.5715
$ws.Range("K138").Value = 6205.7145
$ws.Range("M138").Value = -1065.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5096.25
$ws.Range("I70").Value = 4999
$ws.Range("J70").Value = 5102.7334
$ws.Range("K70").Value = 4999
$ws.Range("L70").Value = 5102.7334
$ws.Range("M70").Value = -4729
$ws.Range("N70").Value = -5642.7334
$ws.Range("H73").Value = 5096.25
$ws.Range("I73").Value = 4999
$ws.Range("J73").Value = 5102.7334
$ws.Range("K73").Value = 4999
$ws.Range("L73").Value = 5102.7334
$ws.Range("M73").Value = -4063
$ws.Range("N73").Value = -6974.7334
$ws.Range("H102").Value = 2742.8286
$ws.Range("I102").Value = 2084.6667
$ws.Range("J102").Value = 3236.45
$ws.Range("K102").Value = 2084.6667
$ws.Range("L102").Value = 3236.45
$ws.Range("M102").Value = -462.6667000000002
$ws.Range("N102").Value = -6480.45
$ws.Range("H122").Value = 2863.7068
$ws.Range("J122").Value = 5292.3887
$ws.Range("L122").Value = 15877.1661
$ws.Range("N122").Value = -20777.1661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4068.4167
$ws.Range("I7").Value = 4138.273
$ws.Range("J7").Value = 3300
$ws.Range("K7").Value = 4138.273
$ws.Range("L7").Value = 3300
$ws.Range("M7").Value = -4026.273
$ws.Range("N7").Value = -3524
$ws.Range("H22").Value = 1532.1538
$ws.Range("I22").Value = 824.1429000000001
$ws.Range("J22").Value = 2358.1667
$ws.Range("K22").Value = 824.1429000000001
$ws.Range("L22").Value = 2358.1667
$ws.Range("M22").Value = -529.1429000000001
$ws.Range("N22").Value = -2948.1667
$ws.Range("H27").Value = 1532.1538
$ws.Range("I27").Value = 824.1429000000001
$ws.Range("J27").Value = 2358.1667
$ws.Range("K27").Value = 824.1429000000001
$ws.Range("L27").Value = 2358.1667
$ws.Range("M27").Value = -717.1429000000001
$ws.Range("N27").Value = -2572.1667
$ws.Range("H40").Value = 9324.532999999999
$ws.Range("I40").Value = 7017.3335
$ws.Range("K40").Value = 7017.3335
$ws.Range("M40").Value = -6881.3335
$ws.Range("H126").Value = 4068.4167
$ws.Range("I126").Value = 4138.273
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 12414.819
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -9944.819
$ws.Range("N126").Value = -14840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2621.196
$ws.Range("I107").Value = 2390.0557
$ws.Range("J107").Value = 3175.9333
$ws.Range("K107").Value = 7170.1671
$ws.Range("L107").Value = 9527.7999
$ws.Range("M107").Value = -5250.1671
$ws.Range("N107").Value = -13367.7999
$ws.Range("H122").Value = 2125.9048
$ws.Range("I122").Value = 2029.5946
$ws.Range("K122").Value = 6088.783799999999
$ws.Range("M122").Value = -3638.783799999999
